$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of election-season data to append (2020 election season).
$newRow = 13

# Mark column B as text *before* writing the value so that the
# numeric-looking string "2020" is stored as text (matching the existing
# Description column entries such as "2018", "2017", etc.) instead of being
# auto-converted to a number by Excel.
$ws.Range("B" + $newRow).NumberFormat = "@"
$ws.Range("B" + $newRow).Value = "2020"

# Copy the formatting (fonts, fills, borders, alignment) of the last
# existing data row down into the new row so the new row matches the look
# of the rest of the table.
$ws.Range("A12:D12").Copy()
$ws.Range(("A" + $newRow + ":D" + $newRow)).PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the actual data values for the new row.
$ws.Range("A" + $newRow).Value = 22
$ws.Range("C" + $newRow).Value = 1
$ws.Range("D" + $newRow).Value = 0
